$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook contains a "Tabela de Requisitos" (requirements table) with columns
# Requisitos (ID) | Descricao | Classificacao, rows 2..23 (22 requirement rows) inside
# an Excel Table (Tabela1) spanning A1:C23.
#
# This change removes two requirement rows that simply described "botao" (button)
# behaviors that were folded/removed from the requirements list:
#   - "Cadastro deve ter um Botao de cadastra-se."   (row 8,  id RF7)
#   - "Login deve ter um Botao de logar-se."          (row 10, id RF9)
# and then renumbers the remaining requirement ids sequentially
# (RF1..RF15, RNF16..RNF20), since the functional (RF) and non functional (RNF)
# requirements share one continuously incrementing counter.

# Delete the higher-numbered row first so the second row index is unaffected.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(8).Delete()

# Renumber column A (Requisitos) sequentially for the remaining 20 requirement rows.
for ($i = 2; $i -le 16; $i++) {
    $ws.Cells.Item($i, 1).Value = "RF" + ($i - 1)
}
for ($i = 17; $i -le 21; $i++) {
    $ws.Cells.Item($i, 1).Value = "RNF" + ($i - 1)
}

# Match the saved selection state recorded in the edited workbook.
$ws.Range("A21").Select()
